$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing milestone date in B5 (6/7/2020 -> 7/6/2020)
$ws.Range("B5").Value = 44018

# Add a new milestone row 6: N=2, Fecha=44049 (8/6/2020), Descripcion = new text
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 44049

# Copy the date number format from B5 onto B6 so it reuses the same style
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C6").Value = "Entendí cómo se ampliaban las opciones en la cadena. Si no encuentro un strike específico, ponerle ALL para que me muestre todo. Ahí ya me va a salir todo. "

# Restore the view selection to match the post-edit state
$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollColumn = 2
